$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.132.07'
$ws.Range('E2').Value = '  +1.10%  '

$ws.Range('D3').Value = '2.214.26'
$ws.Range('E3').Value = '  -0.50%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '230.69'
$ws.Range('D5').Style = "Normal"

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.616'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.91%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '60.67'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.65%  '

$ws.Range('E8').Value = '  +0.02%  '

$ws.Range('E9').Value = '  -0.23%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0902'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +2.31%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.104'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.34%  '

$ws.Range('D12').Value = '2.548.08'
$ws.Range('E12').Value = '  -0.32%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '15.46'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.66%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '22.09'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.80%  '

$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.58'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.82%  '

$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.795'
$ws.Range('D16').Style = "Normal"

$ws.Range('D17').Value = '2.219.24'
$ws.Range('E17').Value = '  -0.16%  '

$ws.Range('D18').Value = '42.130.63'
$ws.Range('E18').Value = '  +1.37%  '

$ws.Range('D19').Value = '0.0₃0934'
$ws.Range('E19').Value = '  +5.21%  '

$ws.Range('E20').Value = '  +2.57%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '72.01'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.73%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '243.80'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.13%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.01%  '

$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.41'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.54%  '

$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.39'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +5.80%  '

$ws.Range('E26').Value = '  +0.20%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '169.43'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.02%  '

$ws.Range('E28').Value = '  -0.24%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '20.26'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.92%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.45'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.77%  '

$ws.Range('E31').Value = '  +2.89%  '

$ws.Range('E32').Value = '  -1.33%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.96'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.42%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.62'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.06%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0651'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +4.46%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.33'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -3.35%  '

$ws.Range('E37').Value = '  -3.15%  '

$ws.Range('E38').Value = '  -1.22%  '

$ws.Range('E39').Value = '  +6.08%  '

$ws.Range('E40').Value = '  +0.16%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.000232'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -2.59%  '

$ws.Range('E42').Value = '  -1.80%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0959'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.33%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.19'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.53%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '96.84'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.07%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.36'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -9.94%  '

$ws.Range('D47').Value = '1.453.44'
$ws.Range('E47').Value = '  -0.72%  '

$ws.Range('B48').Value = 'HuobiToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.75'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.36%  '

$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.07'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.54%  '

$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '15.97'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.22%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.21'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.92%  '
